$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $start = $r.Start
    $end = $r.End - 1
    $r2 = $d.Range($start, $end)
    $r2.Text = $text
    return $r2
}

function Find-ParaIndex([string]$pattern) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2019", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2)

# 2. Insert English translation after the "Objetivos" body paragraph
$idx = Find-ParaIndex("Apresentar ao aluno o conceito de uma organização*")
$para = $d.Paragraphs.Item($idx)
$para.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx + 1)
$textRange = Set-ParaText $newPara "To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues."
$textRange.Font.Italic = $true

# 3. Insert English translation after the "Programa resumido" body paragraph
$idx = Find-ParaIndex("1 - A Administração das Organizações. 2 - O processo administrativo*")
$para = $d.Paragraphs.Item($idx)
$para.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx + 1)
$textRange = Set-ParaText $newPara "The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes"
$textRange.Font.Italic = $true

# 4. Insert English translation after the "Programa" body paragraph
$idx = Find-ParaIndex("1 - A Administração das organizações - definindo a administração*")
$para = $d.Paragraphs.Item($idx)
$para.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx + 1)
$textRange = Set-ParaText $newPara "- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment."
$textRange.Font.Italic = $true

# 5. Replace the Bibliografia content paragraph with the new reference text
$idx = Find-ParaIndex("Gestão de Negócios: Visões*")
$para = $d.Paragraphs.Item($idx)
Set-ParaText $para "LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014."
